$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the StationID for the second Winter deployment (SDWSC -> DWS)
$ws.Range("B3").Value = "DWS"

# Move the active selection to B3, matching the saved selection in the workbook
$ws.Range("B3").Select()
